$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.192668557167053
$ws.Range("B1").Value = 1.374383211135864
$ws.Range("C1").Value = 1.732284545898438
$ws.Range("D1").Value = 3.419886827468872
$ws.Range("E1").Value = 15
